$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.001.27'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').Value = '  +0.74%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.984.94'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').Value = '  +1.09%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.71'
$ws.Range('D5').ClearFormats()

$ws.Range('E6').Value = '  +1.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.35'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').Value = '  +3.17%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '  +2.87%  '

$ws.Range('E10').Value = '  -1.62%  '

$ws.Range('E11').Value = '  +0.46%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.83'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').Value = '  +7.76%  '

$ws.Range('B13').Value = 'Avalanche'

$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.36'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').Value = '  -0.53%  '

$ws.Range('B14').Value = 'Polygon'

$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.850'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '  +2.37%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.271.60'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '  +0.83%  '

$ws.Range('E16').Value = '  +3.93%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.977.98'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').Value = '  +0.15%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.890.72'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').Value = '  +0.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.37'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').Value = '  +0.55%  '

$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.18'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.48'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').Value = '  +0.58%  '

$ws.Range('E24').Value = '  +2.60%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.38'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').Value = '  +0.55%  '

$ws.Range('E26').Value = '  +1.33%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.35'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').Value = '  +0.79%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.02'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').Value = '  +1.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.59'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').Value = '  +0.81%  '

$ws.Range('E30').Value = '  +17.82%  '

$ws.Range('E31').Value = '  +1.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.87'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').Value = '  +2.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0623'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').Value = '  +0.32%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.54'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').Value = '  +5.53%  '

$ws.Range('E35').Value = '  -0.08%  '

$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('E37').Value = '  -1.47%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.56'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').Value = '  -8.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0999'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').Value = '  +1.24%  '

$ws.Range('E41').Value = '  +0.60%  '

$ws.Range('E42').Value = '  +0.44%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0213'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').Value = '  +0.45%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.40'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').Value = '  +1.00%  '

$ws.Range('B45').Value = 'Maker'

$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.378.52'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').Value = '  +1.18%  '

$ws.Range('B46').Value = 'Aave'

$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.33'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').Value = '  +2.64%  '

$ws.Range('E47').Value = '  +0.30%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.28'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').Value = '  +1.60%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.03'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').Value = '  +14.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.41'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').Value = '  +5.41%  '

$ws.Range('E51').Value = '  -0.74%  '
